$d = $word.ActiveDocument
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    $bp = $p.Range.Shading.BackgroundPatternColor
    Write-Output ($idx.ToString() + ": [" + $t.Trim() + "] shd=" + $bp)
}
